# This script re-shuffles the per-row data (Fecha, Variedad, Calidad, Volumen,
# Precio minimo, Precio maximo, Precio promedio ponderado, Unidad de
# comercializacion, Origen, Precio $/Kg) across rows 2-22 according to a
# fixed permutation, as produced by the weekly refresh of the "Hortaliza,
# Mapocho Venta Directa de Santiago - Sandia" sheet. Columns A, B, C, E, F,
# G, Q, R are left untouched (their values are identical for every row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that get shuffled between rows.
$cols = @("D", "H", "I", "J", "K", "L", "M", "N", "O", "P")

# Capture the current ("before") values for rows 2-22 so we can reassign
# them according to the permutation below without clobbering source data
# while we write.
$snapshot = @{}
for ($r = 2; $r -le 22; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowVals
}

# Maps destination row -> source row (values from the source row, as they
# were before this edit, are copied onto the destination row).
$permutation = @{
    2  = 12
    3  = 13
    4  = 14
    5  = 17
    6  = 22
    7  = 21
    8  = 18
    9  = 19
    10 = 6
    11 = 7
    12 = 8
    13 = 9
    14 = 15
    15 = 5
    16 = 20
    17 = 11
    18 = 16
    19 = 2
    20 = 10
    21 = 3
    22 = 4
}

foreach ($destRow in $permutation.Keys) {
    $srcRow = $permutation[$destRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value2 = $srcVals[$c]
    }
}
